$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 64: becomes the former row 65 data (with new B value) ---
$ws.Range("A64").Value = 111998587
$ws.Range("B64").Value = 85313
$ws.Range("E64").Value = 3595
$ws.Range("F64").Value = "Siljansspindling"
$ws.Range("G64").Value = "Cortinarius dalecarlicus"
$ws.Range("H64").Value = "Brandrud"
$ws.Range("Q64").Value = 485439
$ws.Range("R64").Value = 6995893
$ws.Range("Z64").Value = "16:38"
$ws.Range("AB64").Value = "16:38"

# --- Row 65: becomes the former row 64 data (with new B value) ---
$ws.Range("A65").Value = 111998588
$ws.Range("B65").Value = 89089
$ws.Range("E65").Value = 233196
$ws.Range("F65").Value = "Fjällfotad fingersvamp"
$ws.Range("G65").Value = "Ramaria rufescens"
$ws.Range("H65").Value = "(Schaeff.) Corner"
$ws.Range("Q65").Value = 485479
$ws.Range("R65").Value = 6995888
$ws.Range("Z65").Value = "16:31"
$ws.Range("AB65").Value = "16:31"

# --- Row 66: only Taxonsorteringsordning changes ---
$ws.Range("B66").Value = 90803

# --- Row 67: becomes the former row 68 data (with new B value) ---
$ws.Range("A67").Value = 111998584
$ws.Range("B67").Value = 90789
$ws.Range("D67").Value = "VU"
$ws.Range("E67").Value = 150
$ws.Range("F67").Value = "Grangråticka"
$ws.Range("G67").Value = "Boletopsis leucomelaena"
$ws.Range("H67").Value = "(Pers.) Fayod"
$ws.Range("Q67").Value = 485433
$ws.Range("Z67").Value = "16:53"
$ws.Range("AB67").Value = "16:53"

# --- Row 68: becomes the former row 67 data (with new B value) ---
$ws.Range("A68").Value = 111998589
$ws.Range("B68").Value = 89084
$ws.Range("D68").Value = "NT"
$ws.Range("E68").Value = 256756
$ws.Range("F68").Value = "Blek fingersvamp"
$ws.Range("G68").Value = "Ramaria pallida"
$ws.Range("H68").Value = "(Schaeff.) Ricken"
$ws.Range("Q68").Value = 485479
$ws.Range("Z68").Value = "16:20"
$ws.Range("AB68").Value = "16:20"

# --- Row 69: Taxonsorteringsordning changes, public comment removed ---
$ws.Range("B69").Value = 84941
$ws.Range("AC69").ClearContents()
